$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.272.00'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -6.48%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.297.42'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -7.31%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '554.25'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '179.04'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -9.03%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.586'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.24%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.294.69'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -7.05%  '
$ws.Range('E10').Value = '  -11.97%  '
$ws.Range('E11').Value = '  -7.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.12'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -10.76%  '
$ws.Range('E13').Value = '  -9.39%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.847.22'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.51'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -8.20%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '598.59'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -9.92%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '17.96'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '65.316.81'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -6.47%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.116'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.301.86'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -7.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.34'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -9.94%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.898'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -7.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.32'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -6.45%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '102.26'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.87%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.01'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -6.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.95'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -9.97%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.98'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.45%  '
$ws.Range('E28').Value = '  -9.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.30'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -8.91%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.61'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -10.74%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.31'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -9.39%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.85'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -13.98%  '
$ws.Range('E33').Value = '  -8.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '10.97'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -7.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.808.48'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.13%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.104'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.96%  '
$ws.Range('B37').Value = 'Bittensor'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '529.38'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.64%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.999'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '56.04'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -9.38%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.42'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -8.93%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0₃0703'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -15.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.63'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -10.74%  '
$ws.Range('E43').Value = '  -7.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.336'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -10.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '31.63'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -9.45%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.15'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +16.12%  '
$ws.Range('E47').Value = '  -5.28%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0407'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -11.16%  '
$ws.Range('E49').Value = '  -6.28%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.58'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -10.94%  '
$ws.Range('E51').Value = '  +0.03%  '
